# Add class Generic StackCaro
# Updates the "Demo hoàn chỉnh" row (row 12) actual dates and adjusts
# the active selection/view of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 ("Demo hoàn chỉnh"): update actual start date and set actual end date
$ws.Range("H12").Value = "1 tháng 12"
$ws.Range("I12").Value = "6 tháng 12"

# Update the view: scroll so column C is the left-most visible column,
# and move/resize the selection to I12
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("C1").Column
$ws.Range("I12").Select()
